$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'308.76"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.43%"
$ws.Range("E2").Style = "Normal"
$ws.Range("G2").Value = "'22"
$ws.Range("G2").Style = "Normal"
$ws.Range("D3").Value = "'36.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'2.17%"
$ws.Range("E3").Style = "Normal"
$ws.Range("G3").Value = "'22"
$ws.Range("G3").Style = "Normal"
$ws.Range("D4").Value = "'5.077"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'0.19%"
$ws.Range("E4").Style = "Normal"
$ws.Range("G4").Value = "'22"
$ws.Range("G4").Style = "Normal"
$ws.Range("D5").Value = "'0.08239"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'1.15%"
$ws.Range("E5").Style = "Normal"
$ws.Range("G5").Value = "'22"
$ws.Range("G5").Style = "Normal"
$ws.Range("D6").Value = "'2.017"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'1.98%"
$ws.Range("E6").Style = "Normal"
$ws.Range("G6").Value = "'22"
$ws.Range("G6").Style = "Normal"
$ws.Range("D7").Value = "'7.853"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-0.87%"
$ws.Range("E7").Style = "Normal"
$ws.Range("G7").Value = "'22"
$ws.Range("G7").Style = "Normal"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9331"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'0.48%"
$ws.Range("E8").Style = "Normal"
$ws.Range("G8").Value = "'22"
$ws.Range("G8").Style = "Normal"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1525"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'39.44%"
$ws.Range("E9").Style = "Normal"
$ws.Range("G9").Value = "'22"
$ws.Range("G9").Style = "Normal"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1930"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'0.79%"
$ws.Range("E10").Style = "Normal"
$ws.Range("G10").Value = "'22"
$ws.Range("G10").Style = "Normal"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.09060"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-5.22%"
$ws.Range("E11").Style = "Normal"
$ws.Range("G11").Value = "'22"
$ws.Range("G11").Style = "Normal"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03436"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-5.66%"
$ws.Range("E12").Style = "Normal"
$ws.Range("G12").Value = "'22"
$ws.Range("G12").Style = "Normal"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09842"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.99%"
$ws.Range("E13").Style = "Normal"
$ws.Range("G13").Value = "'22"
$ws.Range("G13").Style = "Normal"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001434"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.46%"
$ws.Range("E14").Style = "Normal"
$ws.Range("G14").Value = "'22"
$ws.Range("G14").Style = "Normal"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005717"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-2.13%"
$ws.Range("E15").Style = "Normal"
$ws.Range("G15").Value = "'22"
$ws.Range("G15").Style = "Normal"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.556"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'2.38%"
$ws.Range("E16").Style = "Normal"
$ws.Range("G16").Value = "'22"
$ws.Range("G16").Style = "Normal"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'4.081"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-1.43%"
$ws.Range("E17").Style = "Normal"
$ws.Range("G17").Value = "'22"
$ws.Range("G17").Style = "Normal"
$ws.Range("D18").Value = "'3.048"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'5.23%"
$ws.Range("E18").Style = "Normal"
$ws.Range("G18").Value = "'22"
$ws.Range("G18").Style = "Normal"
$ws.Range("D19").Value = "'0.3436"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'0.47%"
$ws.Range("E19").Style = "Normal"
$ws.Range("G19").Value = "'22"
$ws.Range("G19").Style = "Normal"
$ws.Range("D20").Value = "'0.1288"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-1.71%"
$ws.Range("E20").Style = "Normal"
$ws.Range("G20").Value = "'22"
$ws.Range("G20").Style = "Normal"
$ws.Range("D21").Value = "'5.003"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-2.50%"
$ws.Range("E21").Style = "Normal"
$ws.Range("G21").Value = "'22"
$ws.Range("G21").Style = "Normal"
$ws.Range("D22").Value = "'0.2371"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'7.79%"
$ws.Range("E22").Style = "Normal"
$ws.Range("G22").Value = "'22"
$ws.Range("G22").Style = "Normal"
$ws.Range("D23").Value = "'0.04461"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-1.76%"
$ws.Range("E23").Style = "Normal"
$ws.Range("G23").Value = "'22"
$ws.Range("G23").Style = "Normal"
$ws.Range("B24").Value = "HotbitToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D24").Value = "'0.004895"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'3.44%"
$ws.Range("E24").Style = "Normal"
$ws.Range("G24").Value = "'22"
$ws.Range("G24").Style = "Normal"
$ws.Range("B25").Value = "BitKan"
$ws.Range("C25").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D25").Value = "'0.001194"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-2.90%"
$ws.Range("E25").Style = "Normal"
$ws.Range("G25").Value = "'22"
$ws.Range("G25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001216"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-3.17%"
$ws.Range("E26").Style = "Normal"
$ws.Range("G26").Value = "'22"
$ws.Range("G26").Style = "Normal"
$ws.Range("D27").Value = "'0.0004386"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'-1.83%"
$ws.Range("E27").Style = "Normal"
$ws.Range("G27").Value = "'22"
$ws.Range("G27").Style = "Normal"
$ws.Range("G28").Value = "'22"
$ws.Range("G28").Style = "Normal"
$ws.Range("G29").Value = "'22"
$ws.Range("G29").Style = "Normal"
$ws.Range("G30").Value = "'22"
$ws.Range("G30").Style = "Normal"
$ws.Range("G31").Value = "'22"
$ws.Range("G31").Style = "Normal"
$ws.Range("G32").Value = "'22"
$ws.Range("G32").Style = "Normal"
$ws.Range("G33").Value = "'22"
$ws.Range("G33").Style = "Normal"
$ws.Range("G34").Value = "'22"
$ws.Range("G34").Style = "Normal"
$ws.Range("G35").Value = "'22"
$ws.Range("G35").Style = "Normal"
$ws.Range("G36").Value = "'22"
$ws.Range("G36").Style = "Normal"
$ws.Range("G37").Value = "'22"
$ws.Range("G37").Style = "Normal"
$ws.Range("G38").Value = "'22"
$ws.Range("G38").Style = "Normal"
$ws.Range("D39").Value = "'0.02029"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'4.56%"
$ws.Range("E39").Style = "Normal"
$ws.Range("G39").Value = "'22"
$ws.Range("G39").Style = "Normal"
$ws.Range("D40").Value = "'0.04853"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'0.22%"
$ws.Range("E40").Style = "Normal"
$ws.Range("G40").Value = "'22"
$ws.Range("G40").Style = "Normal"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.007506"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-2.26%"
$ws.Range("E41").Style = "Normal"
$ws.Range("G41").Value = "'22"
$ws.Range("G41").Style = "Normal"
$ws.Range("B42").Value = "Dexo"
$ws.Range("C42").Value = "https://coinranking.com/coin/QkL_pl546+dexo-dexo"
$ws.Range("D42").Value = "'0.01026"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'6.63%"
$ws.Range("E42").Style = "Normal"
$ws.Range("G42").Value = "'22"
$ws.Range("G42").Style = "Normal"
$ws.Range("D43").Value = "'0.1377"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'0.23%"
$ws.Range("E43").Style = "Normal"
$ws.Range("G43").Value = "'22"
$ws.Range("G43").Style = "Normal"
$ws.Range("D44").Value = "'0.002051"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-3.46%"
$ws.Range("E44").Style = "Normal"
$ws.Range("G44").Value = "'22"
$ws.Range("G44").Style = "Normal"
$ws.Range("D45").Value = "'0.01099"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-2.94%"
$ws.Range("E45").Style = "Normal"
$ws.Range("G45").Value = "'22"
$ws.Range("G45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006096"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-6.37%"
$ws.Range("E46").Style = "Normal"
$ws.Range("G46").Value = "'22"
$ws.Range("G46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000741"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-1.59%"
$ws.Range("E47").Style = "Normal"
$ws.Range("G47").Value = "'22"
$ws.Range("G47").Style = "Normal"
$ws.Range("G48").Value = "'22"
$ws.Range("G48").Style = "Normal"
$ws.Range("D49").Value = "'0.001175"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-9.98%"
$ws.Range("E49").Style = "Normal"
$ws.Range("G49").Value = "'22"
$ws.Range("G49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002076"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-1.59%"
$ws.Range("E50").Style = "Normal"
$ws.Range("G50").Value = "'22"
$ws.Range("G50").Style = "Normal"
$ws.Range("D51").Value = "'0.0001977"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-1.59%"
$ws.Range("E51").Style = "Normal"
$ws.Range("G51").Value = "'22"
$ws.Range("G51").Style = "Normal"
